$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to stay as text so numeric-looking values
# (e.g. "0.9999", "236.01") are not auto-converted to numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '25.986.35'
$ws.Range("E2").Value = '  +0.82%  '
$ws.Range("D3").Value = '1.749.37'
$ws.Range("E3").Value = '  -0.10%  '
$ws.Range("D4").Value = '0.9999'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '236.01'
$ws.Range("E5").Value = '  -0.41%  '
$ws.Range("D6").Value = '0.9999'
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("D7").Value = '0.5223'
$ws.Range("E7").Value = '  +3.04%  '
$ws.Range("D8").Value = '0.2850'
$ws.Range("E8").Value = '  +5.48%  '
$ws.Range("D9").Value = '39.31'
$ws.Range("E9").Value = '  -3.91%  '
$ws.Range("D10").Value = '0.06147'
$ws.Range("E10").Value = '  -0.95%  '
$ws.Range("D11").Value = '1.748.54'
$ws.Range("E11").Value = '  -0.17%  '
$ws.Range("D12").Value = '0.07034'
$ws.Range("E12").Value = '  +1.47%  '
$ws.Range("D13").Value = '15.52'
$ws.Range("E13").Value = '  -0.65%  '
$ws.Range("D14").Value = '0.6485'
$ws.Range("D15").Value = '4.538'
$ws.Range("E15").Value = '  +0.94%  '
$ws.Range("D16").Value = '77.59'
$ws.Range("E16").Value = '  -1.30%  '
$ws.Range("D17").Value = '0.9994'
$ws.Range("E17").Value = '  +0.05%  '
$ws.Range("D18").Value = '0.9996'
$ws.Range("E18").Value = '  +0.04%  '
$ws.Range("D19").Value = '25.983.94'
$ws.Range("E19").Value = '  +0.72%  '
$ws.Range("D20").Value = '11.52'
$ws.Range("E20").Value = '  -1.68%  '
$ws.Range("D21").Value = '0.000006640'
$ws.Range("E21").Value = '  -1.54%  '
$ws.Range("D22").Value = '1.976.55'
$ws.Range("E22").Value = '  +0.31%  '
$ws.Range("D23").Value = '4.173'
$ws.Range("E23").Value = '  +2.76%  '
$ws.Range("D24").Value = '8.672'
$ws.Range("E24").Value = '  +4.93%  '
$ws.Range("D25").Value = '5.165'
$ws.Range("E25").Value = '  -0.67%  '
$ws.Range("E26").Value = '  +1.55%  '
$ws.Range("D27").Value = '1.505'
$ws.Range("E27").Value = '  +3.00%  '
$ws.Range("D28").Value = '1.844'
$ws.Range("E28").Value = '  +2.16%  '
$ws.Range("D29").Value = '15.13'
$ws.Range("E29").Value = '  -0.39%  '
$ws.Range("D30").Value = '102.91'
$ws.Range("E30").Value = '  -0.13%  '
$ws.Range("D31").Value = '0.08320'
$ws.Range("E31").Value = '  +0.46%  '
$ws.Range("D32").Value = '3.668'
$ws.Range("E32").Value = '  -2.14%  '
$ws.Range("D33").Value = '3.437'
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("D34").Value = '0.04478'
$ws.Range("E34").Value = '  +1.74%  '
$ws.Range("E35").Value = '  -1.27%  '
$ws.Range("D36").Value = '0.9880'
$ws.Range("E36").Value = '  -2.16%  '
$ws.Range("D37").Value = '0.6123'
$ws.Range("E37").Value = '  +1.10%  '
$ws.Range("D38").Value = '2.687'
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("D39").Value = '0.01597'
$ws.Range("E39").Value = '  +2.14%  '
$ws.Range("D40").Value = '1.957'
$ws.Range("E40").Value = '  -0.65%  '
$ws.Range("D41").Value = '0.9994'
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").Value = '101.18'
$ws.Range("E42").Value = '  -1.07%  '
$ws.Range("D43").Value = '0.3883'
$ws.Range("E43").Value = '  +0.19%  '
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = '5.091'
$ws.Range("E44").Value = '  +4.57%  '
$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").Value = '0.7368'
$ws.Range("E45").Value = '  -2.10%  '
$ws.Range("E46").Value = '  -0.45%  '
$ws.Range("D47").Value = '6.334'
$ws.Range("E47").Value = '  +5.97%  '
$ws.Range("D48").Value = '0.1122'
$ws.Range("E48").Value = '  +3.00%  '
$ws.Range("D49").Value = '53.04'
$ws.Range("E49").Value = '  +0.86%  '
$ws.Range("D50").Value = '30.12'
$ws.Range("E50").Value = '  -0.63%  '
$ws.Range("D51").Value = '7.651'
$ws.Range("E51").Value = '  +2.21%  '

# Restore the original (default/Normal) style on the Price column so
# only cell values changed, matching the source diff.
$priceRange.Style = "Normal"
